# Apply cryptos-list price/volume updates (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.959.47'
$ws.Range('E2').Value = '  +1.33%  '
$ws.Range('D3').Value = '1.769.51'
$ws.Range('E3').Value = '  +0.95%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.80'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4548'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.18%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3530'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.13'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07387'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.00%  '
$ws.Range('E11').Value = '  +1.32%  '
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.72'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.010'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.194'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').Value = '1.771.39'
$ws.Range('E16').Value = '  +0.75%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.63'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.99%  '
$ws.Range('E18').Value = '  +0.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06440'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('E21').Value = '  -1.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.778'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('D23').Value = '27.979.84'
$ws.Range('E23').Value = '  +1.25%  '
$ws.Range('E24').Value = '  +0.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.106'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '160.49'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.20%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').Value = '1.978.78'
$ws.Range('E28').Value = '  +1.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.180'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.35%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.13'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.082'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09256'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.605'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.661'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.85'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.07%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02284'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06130'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.67%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2091'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6260'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.43%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.181'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  -0.85%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.812'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.25'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.02%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.733'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5858'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.85'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.45%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.938'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('E49').Value = '  +0.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06818'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '73.62'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.69%  '
